# Applies the cryptos price/volume refresh described in the commit:
# "Updated cryptos list on Thu Apr  4 16:14:09 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.518.21"
$ws.Range("E2").Value = "  +2.35%  "

# Row 3
$ws.Range("D3").Value = "'3.367.13"
$ws.Range("E3").Value = "  +0.77%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "'591.09"
$ws.Range("E5").Value = "  +6.07%  "

# Row 6
$ws.Range("D6").Value = "'187.11"
$ws.Range("E6").Value = "  -1.21%  "

# Row 7
$ws.Range("D7").Value = "'0.998"

# Row 8
$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  +2.69%  "

# Row 9
$ws.Range("E9").Value = "  +1.54%  "

# Row 10
$ws.Range("E10").Value = "  +1.04%  "

# Row 11
$ws.Range("D11").Value = "'47.45"
$ws.Range("E11").Value = "  +2.49%  "

# Row 12
$ws.Range("D12").Value = "'0.0000276"
$ws.Range("E12").Value = "  +2.45%  "

# Row 13
$ws.Range("D13").Value = "'3.908.22"
$ws.Range("E13").Value = "  +0.97%  "

# Row 14
$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").Value = "'637.24"
$ws.Range("E14").Value = "  +7.35%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'8.61"
$ws.Range("E15").Value = "  +0.71%  "

# Row 16
$ws.Range("D16").Value = "'67.576.69"
$ws.Range("E16").Value = "  +2.38%  "

# Row 17
$ws.Range("E17").Value = "  +1.22%  "

# Row 18
$ws.Range("D18").Value = "'3.365.45"
$ws.Range("E18").Value = "  +0.86%  "

# Row 19
$ws.Range("D19").Value = "'18.05"
$ws.Range("E19").Value = "  +0.62%  "

# Row 20
$ws.Range("D20").Value = "'11.14"
$ws.Range("E20").Value = "  +1.26%  "

# Row 21
$ws.Range("D21").Value = "'0.911"
$ws.Range("E21").Value = "  +1.35%  "

# Row 22
$ws.Range("D22").Value = "'17.96"
$ws.Range("E22").Value = "  -1.69%  "

# Row 23
$ws.Range("D23").Value = "'5.13"
$ws.Range("E23").Value = "  +2.44%  "

# Row 24
$ws.Range("D24").Value = "'99.48"
$ws.Range("E24").Value = "  +0.99%  "

# Row 25
$ws.Range("E25").Value = "  +1.40%  "

# Row 26
$ws.Range("D26").Value = "'2.86"
$ws.Range("E26").Value = "  +4.95%  "

# Row 27
$ws.Range("D27").Value = "'9.71"
$ws.Range("E27").Value = "  +2.65%  "

# Row 28
$ws.Range("D28").Value = "'32.52"
$ws.Range("E28").Value = "  +6.31%  "

# Row 29
$ws.Range("E29").Value = "  +1.84%  "

# Row 30
$ws.Range("D30").Value = "'6.96"
$ws.Range("E30").Value = "  +4.17%  "

# Row 31
$ws.Range("D31").Value = "'604.10"
$ws.Range("E31").Value = "  +3.90%  "

# Row 32
$ws.Range("D32").Value = "'3.79"
$ws.Range("E32").Value = "  -2.56%  "

# Row 33
$ws.Range("D33").Value = "'4.008.68"
$ws.Range("E33").Value = "  +7.88%  "

# Row 34
$ws.Range("E34").Value = "  +1.55%  "

# Row 35
$ws.Range("E35").Value = "  +1.64%  "

# Row 36
$ws.Range("E36").Value = "  -0.12%  "

# Row 37
$ws.Range("D37").Value = "'56.02"
$ws.Range("E37").Value = "  +0.07%  "

# Row 38
$ws.Range("D38").Value = "'2.80"
$ws.Range("E38").Value = "  +6.11%  "

# Row 39
$ws.Range("E39").Value = "  +3.97%  "

# Row 40
$ws.Range("D40").Value = "'33.72"
$ws.Range("E40").Value = "  +0.29%  "

# Row 41
$ws.Range("D41").Value = "'3.23"
$ws.Range("E41").Value = "  +1.97%  "

# Row 42
$ws.Range("D42").Value = "'0.0₃0703"
$ws.Range("E42").Value = "  +0.57%  "

# Row 43
$ws.Range("D43").Value = "'3.41"
$ws.Range("E43").Value = "  -0.24%  "

# Row 44
$ws.Range("E44").Value = "  +1.17%  "

# Row 45
$ws.Range("D45").Value = "'0.0423"
$ws.Range("E45").Value = "  +1.87%  "

# Row 46
$ws.Range("E46").Value = "  +1.34%  "

# Row 47
$ws.Range("E47").Value = "  +1.09%  "

# Row 48
$ws.Range("E48").Value = "  +0.24%  "

# Row 49
$ws.Range("E49").Value = "  +10.77%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'128.49"
$ws.Range("E50").Value = "  +3.08%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'7.73"
$ws.Range("E51").Value = "  +4.18%  "
